$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-29 Wednesday", "2025-01-30 Thursday"),
    @("436×4=", "657×3="),
    @("385×6=", "227×9="),
    @("765×2=", "116×7="),
    @("283×2=", "782×9="),
    @("319×4=", "852×4="),
    @("659×5=", "328×9="),
    @("492×6=", "591×7="),
    @("238×7=", "153×4="),
    @("590×6=", "529×2="),
    @("131×5=", "846×4="),
    @("367×2=", "149×6="),
    @("819×5=", "196×4="),
    @("719×4=", "522×7="),
    @("182×6=", "217×8="),
    @("564×9=", "286×7="),
    @("462×9=", "536×6="),
    @("381×4=", "906×8="),
    @("189×2=", "906×6="),
    @("563×5=", "599×2="),
    @("345×2=", "493×7="),
    @("892×8=", "117×6="),
    @("136×2=", "870×2="),
    @("516×5=", "247×2="),
    @("943×2=", "976×5="),
    @("801×8=", "714×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
